$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title strings (Volume number and week range) ---
$volText = $ws.Range("A8").Value()
$ws.Range("A8").Value = $volText.Replace("17", "18")

$weekText = $ws.Range("C9").Value()
$weekText = $weekText.Replace("4/24/2023", "5/1/2023").Replace("4/30/2023", "5/7/2023")
$ws.Range("C9").Value = $weekText

# --- Update data table cells ---
$ws.Range("C23").Copy($ws.Range("D14"))
$ws.Range("L23").Copy($ws.Range("E14"))
$ws.Range("C23").Copy($ws.Range("G15"))
$ws.Range("L23").Copy($ws.Range("H15"))
$ws.Range("I23").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1
$ws.Range("I23").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 2
$ws.Range("K23").Copy($ws.Range("E16"))
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = -50
$ws.Range("I16").Value = 24
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = -7.692307692307
$ws.Range("L16").Value = 41.176470588235
$ws.Range("M16").Value = -36.842105263157
$ws.Range("N16").Value = -90.438247011952
$ws.Range("I23").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 33
$ws.Range("J17").Value = 34
$ws.Range("K17").Value = -2.941176470588
$ws.Range("L17").Value = 57.142857142857
$ws.Range("M17").Value = 73.684210526315
$ws.Range("N17").Value = -23.255813953488
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 46
$ws.Range("J18").Value = 37
$ws.Range("K18").Value = 24.324324324324
$ws.Range("L18").Value = 53.333333333333
$ws.Range("M18").Value = 6.976744186046
$ws.Range("N18").Value = -90.495867768595
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -63.636363636363
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = -35.897435897435
$ws.Range("I19").Value = 158
$ws.Range("J19").Value = 190
$ws.Range("K19").Value = -16.842105263157
$ws.Range("L19").Value = 25.396825396825
$ws.Range("M19").Value = 22.480620155038
$ws.Range("N19").Value = -52.976190476190
$ws.Range("C20").Value = 2
$ws.Range("C23").Copy($ws.Range("D20"))
$ws.Range("L23").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -44.444444444444
$ws.Range("I20").Value = 37
$ws.Range("K20").Value = 54.166666666666
$ws.Range("L20").Value = 105.555555555556
$ws.Range("M20").Value = -9.756097560975
$ws.Range("N20").Value = -96.967213114754
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -35.294117647058
$ws.Range("G21").Value = 70
$ws.Range("H21").Value = -21.428571428571
$ws.Range("I21").Value = 299
$ws.Range("J21").Value = 318
$ws.Range("K21").Value = -5.974842767295
$ws.Range("L21").Value = 39.069767441860
$ws.Range("M21").Value = 10.332103321033
$ws.Range("N21").Value = -87.211291702309
$ws.Range("I23").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("I23").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K23").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("I23").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 1
$ws.Range("K23").Copy($ws.Range("H22"))
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 10
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = -16.666666666666
$ws.Range("L22").Value = 150
$ws.Range("M22").Value = 25
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = -32.653061224489
$ws.Range("F24").Value = 121
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = -9.022556390977
$ws.Range("I24").Value = 566
$ws.Range("J24").Value = 638
$ws.Range("K24").Value = -11.285266457680
$ws.Range("L24").Value = 21.982758620689
$ws.Range("M24").Value = 72.036474164133
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 35
$ws.Range("I25").Value = 84
$ws.Range("J25").Value = 67
$ws.Range("K25").Value = 25.373134328358
$ws.Range("L25").Value = 47.368421052631
$ws.Range("M25").Value = 16.666666666666
$ws.Range("I23").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("K23").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = -37.5
$ws.Range("L26").Value = -16.666666666666
$ws.Range("I23").Copy($ws.Range("C27"))
$ws.Range("C27").Value = 1
$ws.Range("C23").Copy($ws.Range("D27"))
$ws.Range("L23").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 9
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -35.714285714285
$ws.Range("C23").Copy($ws.Range("D28"))
$ws.Range("L23").Copy($ws.Range("E28"))
$ws.Range("C23").Copy($ws.Range("D29"))
$ws.Range("L23").Copy($ws.Range("E29"))
$ws.Range("C23").Copy($ws.Range("C30"))
